$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1963
$ws.Range("I2").Value = 5317
$ws.Range("J2").Value = 21898
$ws.Range("K2").Value = 104
$ws.Range("L2").Value = 6058
$ws.Range("M2").Value = 364
$ws.Range("N2").Value = 3781
$ws.Range("O2").Value = 14
$ws.Range("P2").Value = 84
$ws.Range("Q2").Value = 33
$ws.Range("R2").Value = 272
$ws.Range("S2").Value = 2382
$ws.Range("T2").Value = 3902
$ws.Range("U2").Value = 298
$ws.Range("V2").Value = 33854
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 34323
$ws.Range("Y2").Value = 50
$ws.Range("Z2").Value = 530
$ws.Range("AA2").Value = 226
